# Add a new "2023" column (Q) to the water-loss-during-transportation sheet,
# mirroring the formatting of the existing "2022" column (P), and tighten up
# row heights / column width as in the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- New column Q: header year + one data value per row -------------------
# Each row's Q cell should carry the same style as the corresponding P cell
# in that row, so copy P's formatting over to Q first, then write the value.

$values = @{
    3  = 2023
    5  = 2385.9
    6  = 112.1
    7  = 267.89999999999998
    8  = 230.9
    9  = 249.7
    10 = 287
    11 = 334.7
    12 = 851
    13 = 48.5
    14 = 4.2
    16 = 26.890545708088244
    17 = 15.490056759274875
    18 = 22.218388220841799
    19 = 29.614327895683314
    20 = 30.104452089276922
    21 = 21.825966598728439
    22 = 32.351574864874735
    23 = 30.810022297218843
    24 = 29.193884213235311
    25 = 7.4362892319581295
}

# Row 15 has no value in P (just formatting, style 23) -> Q15 mirrors that:
# formatted but empty. Row 4 has no P cell at all, so it is skipped entirely.
$formatOnlyRows = @(15)

foreach ($row in ($values.Keys + $formatOnlyRows)) {
    $srcCell = $ws.Cells.Item($row, 16)   # column P
    $dstCell = $ws.Cells.Item($row, 17)   # column Q

    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial($xlPasteFormats) | Out-Null

    if ($values.ContainsKey($row)) {
        $dstCell.Value = $values[$row]
    }
}

# --- Row heights: rows 4-25 become an explicit 15pt custom height ---------
foreach ($row in 4..25) {
    $ws.Rows.Item($row).RowHeight = 15
}

# --- Column A:C width tightened slightly -----------------------------------
$ws.Columns("A:C").ColumnWidth = 35.5

# --- Reset the view selection back to the top-left cell -------------------
$ws.Range("A1").Select() | Out-Null
